$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update "last updated" timestamp text (A1)
$ws.Range("A1").Value = "Datos actualizados a 26 de Mayo de 2020 a las 23:35"

# Row 4: Estados Unidos - refreshed totals
$ws.Range("B4").Value = 1722384
$ws.Range("C4").Value = 16158
$ws.Range("D4").Value = 474923
$ws.Range("E4").Value = 1147002
$ws.Range("G4").Value = 654
$ws.Range("H4").Value = 100459

# Row 16: Canada - refreshed totals
$ws.Range("D16").Value = 45250
$ws.Range("E16").Value = 34726
$ws.Range("G16").Value = 93
$ws.Range("H16").Value = 6638

# Row 88: Gabon - refreshed totals
$ws.Range("B88").Value = 2238
$ws.Range("C88").Value = 103
$ws.Range("D88").Value = 593
$ws.Range("E88").Value = 1631

# Rows 153/154: Uganda overtakes Yemen in ranking, so the two countries swap
# rows; Uganda's figures are refreshed while Yemen keeps its prior figures.
$ws.Range("A153").Value = "Uganda"
$ws.Range("B153").Value = 253
$ws.Range("C153").Value = 31
$ws.Range("D153").Value = 69
$ws.Range("E153").Value = 184
$ws.Range("G153").Value = 0
$ws.Range("H153").Value = 0

$ws.Range("A154").Value = "Yemen"
$ws.Range("B154").Value = 249
$ws.Range("C154").Value = 16
$ws.Range("D154").Value = 10
$ws.Range("E154").Value = 190
$ws.Range("G154").Value = 5
$ws.Range("H154").Value = 49
